$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give every data row (1 through 34) a custom height of 22.5pt
$ws.Range("A1:E34").RowHeight = 22.5

# Fill in the previously-empty "Sprint Number" (column E) values for
# requirements 6.0, 6.1, 7.0 and 7.1 (rows 17, 18, 20, 21)
$ws.Range("E17").Value = 3
$ws.Range("E18").Value = 3
$ws.Range("E20").Value = 3
$ws.Range("E21").Value = 3

# Apply the same formatting already used by the other Sprint Number cells
# (e.g. E11) to the newly populated cells, reusing the existing cell style
$ws.Range("E11").Copy()
$ws.Range("E17").PasteSpecial(-4122)
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E21").PasteSpecial(-4122)

$excel.CutCopyMode = 0
